$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.02406370639801
$ws.Range("B1").Value = 2.424098014831543
$ws.Range("C1").Value = 5.205649852752686
$ws.Range("D1").Value = 2.270311594009399
$ws.Range("E1").Value = 1.337273120880127
